# Refresh the cryptocurrency price/volume snapshot values for the latest
# GitHub Actions scrape run (Wed Feb 15 13:50:03 UTC 2023).
#
# Price (column D) and Volume(1h) (column E) are stored as plain text in this
# sheet (e.g. "303.56", "3.19%"), not as numbers/percentages. Assigning a
# numeric-looking string straight to Range.Value would make Excel silently
# reinterpret it as a real number, so each cell is briefly forced to the
# "@" (Text) number format while the new value is written, then restored to
# the workbook's normal (default) style so formatting is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $text) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "303.56"
Set-TextValue "E2" "3.19%"
Set-TextValue "D3" "42.89"
Set-TextValue "E3" "7.06%"
Set-TextValue "D4" "5.039"
Set-TextValue "E4" "0.46%"
Set-TextValue "D5" "0.07684"
Set-TextValue "E5" "4.40%"
Set-TextValue "D6" "4.405"
Set-TextValue "E6" "2.53%"
Set-TextValue "E7" "4.67%"
Set-TextValue "D8" "1.058"
Set-TextValue "E8" "14.53%"
Set-TextValue "D10" "0.1232"
Set-TextValue "E10" "4.43%"
Set-TextValue "D11" "0.1851"
Set-TextValue "E11" "3.75%"
Set-TextValue "D12" "0.09010"
Set-TextValue "E12" "3.72%"
Set-TextValue "D13" "0.04164"
Set-TextValue "E13" "-1.71%"
Set-TextValue "D14" "0.1045"
Set-TextValue "E14" "-0.85%"
Set-TextValue "D15" "0.001269"
Set-TextValue "E15" "1.23%"
Set-TextValue "D16" "0.005766"
Set-TextValue "E16" "-3.58%"
Set-TextValue "E17" "1,892.27%"
Set-TextValue "D18" "3.325"
Set-TextValue "E18" "-1.25%"
Set-TextValue "D19" "0.3341"
Set-TextValue "D20" "8.402"
Set-TextValue "E20" "6.60%"
Set-TextValue "D21" "0.1404"
Set-TextValue "E21" "1.58%"
Set-TextValue "D22" "0.2892"
Set-TextValue "E22" "2.86%"
Set-TextValue "D23" "0.04151"
Set-TextValue "E23" "5.35%"
Set-TextValue "D24" "0.001273"
Set-TextValue "E24" "0.31%"
Set-TextValue "D25" "0.004501"
Set-TextValue "E25" "17.97%"
Set-TextValue "D26" "0.0001346"
Set-TextValue "E26" "9.22%"
Set-TextValue "D38" "0.02455"
Set-TextValue "E38" "4.88%"
Set-TextValue "D39" "0.05274"
Set-TextValue "E39" "3.97%"
Set-TextValue "D40" "0.005930"
Set-TextValue "E40" "-3.49%"
Set-TextValue "D41" "0.007655"
Set-TextValue "E41" "-1.65%"
Set-TextValue "E42" "4.51%"
Set-TextValue "D43" "0.007346"
Set-TextValue "E43" "-0.56%"
Set-TextValue "D44" "0.008361"
Set-TextValue "E44" "16.36%"
Set-TextValue "D45" "0.3020"
Set-TextValue "E45" "3.24%"
Set-TextValue "D46" "0.00006632"
Set-TextValue "E46" "8.22%"
Set-TextValue "D47" "0.00000000748"
Set-TextValue "E47" "-0.48%"
Set-TextValue "E48" "-0.19%"
Set-TextValue "E49" "-17.78%"
Set-TextValue "D50" "0.00002094"
Set-TextValue "E50" "-0.48%"
Set-TextValue "D51" "0.0001995"
Set-TextValue "E51" "-0.48%"
